# Regenerate the "K" column (G) of the save-data sheet with freshly
# calculated strikeout values (K), replacing the old Strike# values.
# This mirrors: "regen save_data to use K instead of Strike#, regen
# std/mean, calc and write s_vals".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New calculated K values for rows 2..63 (column G), in row order.
$newK = @(1,1,0,1,1,1,2,0,0,0,2,0,0,2,2,0,0,0,0,1,1,0,1,0,3,0,0,1,0,0,0,0,1,1,0,0,0,1,0,1,0,1,1,0,0,0,0,1,0,1,0,1,1,1,0,0,0,0,1,2,1,1)

$startRow = 2
for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $newK[$i]
}
